$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.005.54'
$ws.Range("E2").Value = '  -0.07%  '

$ws.Range("D3").Value = '3.935.18'
$ws.Range("E3").Value = '  +3.62%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Formula = '="604.43"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E5").Value = '  +0.73%  '

$ws.Range("D6").Formula = '="168.88"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E6").Value = '  +2.85%  '

$ws.Range("D7").Value = '3.933.41'
$ws.Range("E7").Value = '  +3.67%  '

$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("E9").Value = '  -0.02%  '

$ws.Range("D10").Formula = '="0.170"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E10").Value = '  +1.17%  '

$ws.Range("D11").Formula = '="6.50"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E11").Value = '  +3.08%  '

$ws.Range("D12").Formula = '="0.467"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E12").Value = '  +1.84%  '

$ws.Range("D13").Formula = '="0.0000256"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E13").Value = '  +4.55%  '

$ws.Range("D14").Formula = '="37.71"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E14").Value = '  +1.90%  '

$ws.Range("D15").Value = '4.589.14'
$ws.Range("E15").Value = '  +3.55%  '

$ws.Range("D16").Value = '3.908.82'
$ws.Range("E16").Value = '  +4.35%  '

$ws.Range("D17").Value = '69.008.95'
$ws.Range("E17").Value = '  -0.25%  '

$ws.Range("E18").Value = '  +0.40%  '

$ws.Range("D19").Formula = '="17.41"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E19").Value = '  +1.22%  '

$ws.Range("E20").Value = '  -1.71%  '

$ws.Range("D21").Formula = '="11.01"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E21").Value = '  -2.70%  '

$ws.Range("D22").Formula = '="494.58"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E22").Value = '  +1.62%  '

$ws.Range("E23").Value = '  +1.93%  '

$ws.Range("D24").Formula = '="0.0000167"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E24").Value = '  +6.09%  '

$ws.Range("D25").Formula = '="84.93"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E25").Value = '  +0.62%  '

$ws.Range("D26").Formula = '="2.28"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E26").Value = '  +1.33%  '

$ws.Range("E27").Value = '  +0.57%  '

$ws.Range("D28").Formula = '="10.26"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E28").Value = '  +2.36%  '

$ws.Range("E29").Value = '  +0.11%  '

$ws.Range("E30").Value = '  +0.90%  '

$ws.Range("D31").Value = '4.083.59'
$ws.Range("E31").Value = '  +3.26%  '

$ws.Range("E32").Value = '  +0.67%  '

$ws.Range("D33").Formula = '="7.84"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E33").Value = '  -1.96%  '

$ws.Range("D34").Formula = '="32.17"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E34").Value = '  +1.86%  '

$ws.Range("D35").Value = '3.892.48'
$ws.Range("E35").Value = '  +4.05%  '

$ws.Range("D36").Formula = '="0.108"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E36").Value = '  +0.85%  '

$ws.Range("E37").Value = '  +1.69%  '

$ws.Range("D38").Formula = '="6.01"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E38").Value = '  +2.80%  '

$ws.Range("E39").Value = '  +0.48%  '

$ws.Range("D40").Formula = '="3.32"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E40").Value = '  +9.60%  '

$ws.Range("E41").Value = '  -0.02%  '

$ws.Range("D42").Formula = '="0.323"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E42").Value = '  +1.50%  '

$ws.Range("D43").Formula = '="441.31"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E43").Value = '  +0.95%  '

$ws.Range("E44").Value = '  +1.36%  '

$ws.Range("D45").Formula = '="48.15"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("E46").Value = '  +2.99%  '

$ws.Range("B48").Value = 'FLOKI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D48").Formula = '="0.000273"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E48").Value = '  +20.52%  '

$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").Formula = '="143.17"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E49").Value = '  +1.16%  '

$ws.Range("D50").Value = '2.824.45'
$ws.Range("E50").Value = '  +0.27%  '

$ws.Range("D51").Formula = '="0.0359"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E51").Value = '  +1.81%  '
